# Apply the IPS.HAJJ.CONSENT ValueSet update:
#  - rename the "Include from IPS.HAJJ.CONSENT" sheet to "Include #0"
#  - bump Version / Date on the Metadata sheet
#  - add a new "Jurisdiction" metadata row
#  - replace the per-concept "Concept" table on the Include sheet with a
#    single "Codes" / "All codes" summary (keeping the System URI row)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Include from IPS.HAJJ.CONSENT"

# --- workbook-level: rename the include sheet -----------------------------
$ws2.Name = "Include #0"

# --- Metadata sheet ---------------------------------------------------------
# Version: 1.0.0 -> 2.0.2
$ws1.Range("B3").Value = "2.0.2"

# Date: 2024-04-25T17:24:48+00:00 -> 2025-02-13T16:11:24+00:00
$ws1.Range("B8").Value = "2025-02-13T16:11:24+00:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws1.Rows("11:11").Insert()
$ws1.Range("A10:B10").Copy($ws1.Range("A11:B11"))
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- Include sheet: collapse the per-concept list into a code summary ------
# Original layout:
#   1 Concept      | Description
#   2 CONSENT-KSA  | Consent for Kingdom of Saudi Arabia
#   3 CONSENT-IDN  | Consent for Indonesia
#   4 CONSENT-MYS  | Consent for Malaysia
#   5 CONSENT-OMN  | Consent for Oman
#   6 (blank)      | (blank)
#   7 System URI   | http://smart.who.int/ips-pilgrimage/CodeSystem/IPS.HAJJ.CONSENT
#
# Target layout:
#   1 Codes
#   2 All codes
#   3 (blank)      | (blank)
#   4 System URI   | http://smart.who.int/ips-pilgrimage/CodeSystem/IPS.HAJJ.CONSENT

# Drop column B entirely first so the two new header rows end up with only
# a column-A cell (no B1/B2 left behind).
$ws2.Columns("B").Delete()

# Remove the three per-country concept rows (old rows 3-5); the blank
# separator row and the System URI row (old rows 6-7) shift up to become
# rows 3-4.
$ws2.Rows("3:5").Delete()

# Re-create column B only where the target still needs it (rows 3 and 4).
$ws2.Range("A3").Copy($ws2.Range("B3"))
$ws2.Range("A4").Copy($ws2.Range("B4"))
$ws2.Range("B4").Value = "http://smart.who.int/ips-pilgrimage/CodeSystem/IPS.HAJJ.CONSENT"

# Rewrite the two header rows as "Codes" / "All codes".
$ws2.Range("A1").Value = "Codes"
$ws2.Range("A2").Value = "All codes"

Write-Output "IPS.HAJJ.CONSENT metadata + include sheet updated"
